$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("G1").Select()
$win.FreezePanes = $true
